# "finestra incidenza 7gg centrata su ultimo g"
# Shift the rolling 7-day window (col C = 7-day sum of new cases, col D = incidence
# per 100k inhabitants) so it is a trailing window ending on the row's own date
# (i.e. "centered on the last day") instead of a window centered on the row's date.
#
# Population of the comune (used for the per-100k-inhabitants figure) is derived
# from the existing data: C5=9 -> D5=37.31033910952657 => population = 24122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$population = 24122
$firstDataRow = 2
$lastDataRow = 184
$windowSize = 7

# Read every "nuovi pos." (col B) value up front.
$newCases = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $newCases[$r] = $ws.Cells.Item($r, 2).Value2
}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $windowStart = $r - $windowSize + 1
    if ($windowStart -lt $firstDataRow) {
        # Not enough prior days to fill a full trailing 7-day window: leave blank.
        # Only touch cells that actually currently hold a number (rows that used
        # to show a centered-window value but can no longer be computed); cells
        # that were already blank stay untouched.
        if ($ws.Cells.Item($r, 3).Value2 -ne "") {
            $ws.Cells.Item($r, 3).ClearContents()
        }
        if ($ws.Cells.Item($r, 4).Value2 -ne "") {
            $ws.Cells.Item($r, 4).ClearContents()
        }
    } else {
        $sum = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $sum = $sum + $newCases[$i]
        }
        $ws.Cells.Item($r, 3).Value = $sum
        $ws.Cells.Item($r, 4).Value = $sum * 100000 / $population
    }
}
